$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "1.00", "0.0000241", "3.320.80") are preserved verbatim as text,
# matching the source data which stores these as plain strings.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '92.459.10'
$ws.Range("E2").Value = '  -6.00%  '
$ws.Range("D3").Value = '3.320.80'
$ws.Range("E3").Value = '  -5.15%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '228.56'
$ws.Range("E5").Value = '  -9.97%  '
$ws.Range("D6").Value = '623.29'
$ws.Range("E6").Value = '  -6.67%  '
$ws.Range("D7").Value = '1.34'
$ws.Range("E7").Value = '  -10.21%  '
$ws.Range("D8").Value = '0.380'
$ws.Range("E8").Value = '  -11.62%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '0.919'
$ws.Range("E10").Value = '  -13.18%  '
$ws.Range("D11").Value = '3.318.21'
$ws.Range("E11").Value = '  -5.16%  '
$ws.Range("D12").Value = '0.192'
$ws.Range("E12").Value = '  -8.88%  '
$ws.Range("D13").Value = '39.62'
$ws.Range("E13").Value = '  -13.52%  '
$ws.Range("B14").Value = 'WrappedBTC'
$ws.Range("C14").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D14").Value = '92.347.77'
$ws.Range("E14").Value = '  -5.94%  '
$ws.Range("B15").Value = 'Toncoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D15").Value = '5.86'
$ws.Range("E15").Value = '  -6.28%  '
$ws.Range("D16").Value = '3.940.21'
$ws.Range("E16").Value = '  -5.33%  '
$ws.Range("D17").Value = '0.0000241'
$ws.Range("E17").Value = '  -8.50%  '
$ws.Range("D18").Value = '7.83'
$ws.Range("E18").Value = '  -13.14%  '
$ws.Range("D19").Value = '3.324.19'
$ws.Range("E19").Value = '  -5.13%  '
$ws.Range("D20").Value = '16.54'
$ws.Range("E20").Value = '  -11.83%  '
$ws.Range("D21").Value = '10.80'
$ws.Range("E21").Value = '  -8.72%  '
$ws.Range("D22").Value = '486.33'
$ws.Range("E22").Value = '  -7.04%  '
$ws.Range("D23").Value = '0.442'
$ws.Range("E23").Value = '  -16.39%  '
$ws.Range("E24").Value = '  -10.31%  '
$ws.Range("D25").Value = '0.0000182'
$ws.Range("E25").Value = '  -11.20%  '
$ws.Range("D26").Value = '6.16'
$ws.Range("E26").Value = '  -9.53%  '
$ws.Range("D27").Value = '88.77'
$ws.Range("E27").Value = '  -9.81%  '
$ws.Range("D28").Value = '3.520.60'
$ws.Range("E28").Value = '  -4.59%  '
$ws.Range("D29").Value = '11.28'
$ws.Range("E29").Value = '  -11.38%  '
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("D31").Value = '11.09'
$ws.Range("E31").Value = '  -12.21%  '
$ws.Range("D32").Value = '2.61'
$ws.Range("E32").Value = '  -9.50%  '
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.55%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.130'
$ws.Range("E34").Value = '  -11.69%  '
$ws.Range("D35").Value = '0.168'
$ws.Range("E35").Value = '  -12.37%  '
$ws.Range("D36").Value = '28.10'
$ws.Range("E36").Value = '  -9.45%  '
$ws.Range("D37").Value = '0.517'
$ws.Range("E37").Value = '  -14.15%  '
$ws.Range("B38").Value = 'USDe'
$ws.Range("C38").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").Value = '7.32'
$ws.Range("E39").Value = '  -8.86%  '
$ws.Range("D40").Value = '511.60'
$ws.Range("E40").Value = '  -3.51%  '
$ws.Range("D41").Value = '1.37'
$ws.Range("E41").Value = '  -10.82%  '
$ws.Range("D42").Value = '0.145'
$ws.Range("E42").Value = '  -6.87%  '
$ws.Range("D43").Value = '0.861'
$ws.Range("E43").Value = '  -5.57%  '
$ws.Range("D44").Value = '23.99'
$ws.Range("E44").Value = '  -1.84%  '
$ws.Range("B45").Value = 'MantraDAO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D45").Value = '3.54'
$ws.Range("E45").Value = '  -3.17%  '
$ws.Range("B46").Value = 'ImmutableX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D46").Value = '1.64'
$ws.Range("E46").Value = '  -8.06%  '
$ws.Range("D47").Value = '5.36'
$ws.Range("E47").Value = '  -7.97%  '
$ws.Range("D48").Value = '2.11'
$ws.Range("E48").Value = '  -4.95%  '
$ws.Range("D49").Value = '0.0386'
$ws.Range("E49").Value = '  -11.69%  '
$ws.Range("D50").Value = '51.73'
$ws.Range("E50").Value = '  -6.49%  '
$ws.Range("E51").Value = '  -5.45%  '

# Restore default (unstyled) formatting on the Price/Volume columns so the
# temporary text-number-format above does not leave a residual style on cells.
$ws.Range("D2:E51").Style = "Normal"
